$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the existing row 542,
# pushing every subsequent row (542-584) down by one (543-585).
$ws.Rows.Item(542).Insert()

# Fill in the new row's data.
$ws.Range("A542").Value = 9
$ws.Range("B542").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C542").Value = "Metropolitana"
$ws.Range("D542").Value = 45106
$ws.Range("E542").Value = 13
$ws.Range("F542").Value = 100112039
$ws.Range("G542").Value = "Ciboulette"
$ws.Range("H542").Value = "Sin especificar"
$ws.Range("I542").Value = "Primera"
$ws.Range("J542").Value = 250
$ws.Range("K542").Value = 1500
$ws.Range("L542").Value = 1700
$ws.Range("M542").Value = 1600
$ws.Range("N542").Value = "$/docena de atados"
$ws.Range("O542").Value = "Región Metropolitana"
$ws.Range("P542").Value = 533
$ws.Range("Q542").Value = 3
$ws.Range("R542").Value = "Hortaliza"
